$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ordering of model names for rows 2..26 (column A)
$names = @(
    "model_4_8_0",
    "model_4_8_22",
    "model_4_8_21",
    "model_4_8_20",
    "model_4_8_19",
    "model_4_8_18",
    "model_4_8_17",
    "model_4_8_16",
    "model_4_8_15",
    "model_4_8_14",
    "model_4_8_13",
    "model_4_8_23",
    "model_4_8_12",
    "model_4_8_10",
    "model_4_8_9",
    "model_4_8_8",
    "model_4_8_7",
    "model_4_8_6",
    "model_4_8_5",
    "model_4_8_4",
    "model_4_8_3",
    "model_4_8_2",
    "model_4_8_1",
    "model_4_8_11",
    "model_4_8_24"
)

# Common numeric values now shared by every data row (B..Q)
$values = @(
    0.445980542904741,
    -0.05067851974480586,
    -0.6284804215699207,
    0.994007028627763,
    0.4605777367007629,
    0.3288898821261522,
    0.6237281562693088,
    0.4122475988555908,
    0.003008683658732021,
    0.2076281412571614,
    0.279043433096329,
    0.5734892170966707,
    0.05025235926527027,
    0.5979038403548141,
    22.22406457805329,
    34.41282282673529
)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $names[$i]
    for ($j = 0; $j -lt $values.Length; $j++) {
        $col = $j + 2
        $ws.Cells.Item($row, $col).Value = $values[$j]
    }
}
